$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- CIN2+ block: prefix generic SE/LL95/UL95 headers with "CIN2+ "
$ws.Range("U1").Value  = "CIN2+ SE immediate"
$ws.Range("V1").Value  = "CIN2+ LL95 immediate"
$ws.Range("W1").Value  = "CIN2+ UL95 immediate"
$ws.Range("Y1").Value  = "CIN2+ SE 1-year"
$ws.Range("Z1").Value  = "CIN2+ LL95 1-year"
$ws.Range("AA1").Value = "CIN2+ UL95 1-year"
$ws.Range("AC1").Value = "CIN2+ SE 2-year"
$ws.Range("AD1").Value = "CIN2+ LL95 2-year"
$ws.Range("AE1").Value = "CIN2+ UL95 2-year"
$ws.Range("AG1").Value = "CIN2+ SE 3-year"
$ws.Range("AH1").Value = "CIN2+ LL95 3-year"
$ws.Range("AI1").Value = "CIN2+ UL95 3-year"
$ws.Range("AK1").Value = "CIN2+ SE 4-year"
$ws.Range("AL1").Value = "CIN2+ LL95 4-year"
$ws.Range("AM1").Value = "CIN2+ UL95 4-year"
$ws.Range("AO1").Value = "CIN2+ SE 5-year"
$ws.Range("AP1").Value = "CIN2+ LL95 5-year"
$ws.Range("AQ1").Value = "CIN2+ UL95 5-year"

# --- CIN3+ block
$ws.Range("AS1").Value = "CIN3+ SE immediate"
$ws.Range("AT1").Value = "CIN3+ LL95 immediate"
$ws.Range("AU1").Value = "CIN3+ UL95 immediate"
$ws.Range("AW1").Value = "CIN3+ SE 1-year"
$ws.Range("AX1").Value = "CIN3+ LL95 1-year"
$ws.Range("AY1").Value = "CIN3+ UL95 1-year"
$ws.Range("BA1").Value = "CIN3+ SE 2-year"
$ws.Range("BB1").Value = "CIN3+ LL95 2-year"
$ws.Range("BC1").Value = "CIN3+ UL95 2-year"
$ws.Range("BE1").Value = "CIN3+ SE 3-year"
$ws.Range("BF1").Value = "CIN3+ LL95 3-year"
$ws.Range("BG1").Value = "CIN3+ UL95 3-year"
$ws.Range("BI1").Value = "CIN3+ SE 4-year"
$ws.Range("BJ1").Value = "CIN3+ LL95 4-year"
$ws.Range("BK1").Value = "CIN3+ UL95 4-year"
$ws.Range("BM1").Value = "CIN3+ SE 5-year"
$ws.Range("BN1").Value = "CIN3+ LL95 5-year"
$ws.Range("BO1").Value = "CIN3+ UL95 5-year"

# --- CANCER block
$ws.Range("BQ1").Value = "CANCER SE immediate"
$ws.Range("BR1").Value = "CANCER LL95 immediate"
$ws.Range("BS1").Value = "CANCER UL95 immediate"
$ws.Range("BU1").Value = "CANCER SE 1-year"
$ws.Range("BV1").Value = "CANCER LL95 1-year"
$ws.Range("BW1").Value = "CANCER UL95 1-year"
$ws.Range("BY1").Value = "CANCER SE 2-year"
$ws.Range("BZ1").Value = "CANCER LL95 2-year"
$ws.Range("CA1").Value = "CANCER UL95 2-year"
$ws.Range("CC1").Value = "CANCER SE 3-year"
$ws.Range("CD1").Value = "CANCER LL95 3-year"
$ws.Range("CE1").Value = "CANCER UL95 3-year"
$ws.Range("CG1").Value = "CANCER SE 4-year"
$ws.Range("CH1").Value = "CANCER LL95 4-year"
$ws.Range("CI1").Value = "CANCER UL95 4-year"
$ws.Range("CK1").Value = "CANCER SE 5-year"
$ws.Range("CL1").Value = "CANCER LL95 5-year"
$ws.Range("CM1").Value = "CANCER UL95 5-year"

# --- Unweighted block: prefix with "Unweighted "
$ws.Range("CR1").Value = "Unweighted %"
$ws.Range("CT1").Value = "Unweighted CIN2+ Prevalence Cases"
$ws.Range("CU1").Value = "Unweighted CIN2+ Incidence Cases"
$ws.Range("CV1").Value = "Unweighted CIN2+ Unknown Cases"
$ws.Range("CW1").Value = "Unweighted Number of CIN3+ Cases"
$ws.Range("CX1").Value = "Unweighted CIN3+ Prevalence Cases"
$ws.Range("CY1").Value = "Unweighted CIN3+ Incidence Cases"
$ws.Range("CZ1").Value = "Unweighted CIN3+ Unknown Cases"
$ws.Range("DA1").Value = "Unweighted Number of Cancer Cases"
$ws.Range("DB1").Value = "Unweighted Cancer Prevalence Cases"
$ws.Range("DC1").Value = "Unweighted Cancer Incidence Cases"
$ws.Range("DD1").Value = "Unweighted Cancer Unknown Cases"

# --- Update the active selection to the header row range (A1:DD1)
$ws.Range("A1:DD1").Select()
